$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.01860366666666667
$ws.Range("H2").Value = 0.055811
$ws.Range("I2").Value = 0.01426592996042112
$ws.Range("J2").Value = 0.0209748055045382
$ws.Range("M2").Value = 98.170451
$ws.Range("N2").Value = 294.511353
$ws.Range("O2").Value = 0.4110278868558457
$ws.Range("P2").Value = 0.4243928215400083
$ws.Range("Q2").Value = 1.826330346920333
$ws.Range("R2").Value = 16.436973122283
$ws.Range("S2").Value = 0.005863695045665394
$ws.Range("T2").Value = 0.008901556889323863

$ws.Range("G3").Value = 0.01860366666666667
$ws.Range("H3").Value = 0.055811
$ws.Range("I3").Value = 0.01426592996042112
$ws.Range("J3").Value = 0.0209748055045382
$ws.Range("O3").Value = 0.1453748421312515
$ws.Range("P3").Value = 0.150101833491052
$ws.Range("Q3").Value = 0.6459476214473334
$ws.Range("R3").Value = 5.813528593026001
$ws.Range("S3").Value = 0.002073907315851711
$ws.Range("T3").Value = 0.003148356763349393

$ws.Range("G4").Value = 0.01860366666666667
$ws.Range("H4").Value = 0.055811
$ws.Range("I4").Value = 0.01426592996042112
$ws.Range("J4").Value = 0.0209748055045382
$ws.Range("M4").Value = 38.82199566666667
$ws.Range("N4").Value = 116.465987
$ws.Range("O4").Value = 0.1625430328561575
$ws.Range("P4").Value = 0.1678282630971171
$ws.Range("Q4").Value = 0.7222314667174444
$ws.Range("R4").Value = 6.500083200456999
$ws.Range("S4").Value = 0.002318827522280372
$ws.Range("T4").Value = 0.003520165176626497

$ws.Range("G5").Value = 0.01860366666666667
$ws.Range("H5").Value = 0.055811
$ws.Range("I5").Value = 0.01426592996042112
$ws.Range("J5").Value = 0.0209748055045382
$ws.Range("M5").Value = 22.5647
$ws.Range("N5").Value = 45.1294
$ws.Range("O5").Value = 0.09447568860141126
$ws.Range("P5").Value = 0.06503176602637677
$ws.Range("Q5").Value = 0.4197861572333334
$ws.Range("R5").Value = 2.5187169434
$ws.Range("S5").Value = 0.001347783556550289
$ws.Range("T5").Value = 0.001364028644019887

$ws.Range("G6").Value = 0.01860366666666667
$ws.Range("H6").Value = 0.055811
$ws.Range("I6").Value = 0.01426592996042112
$ws.Range("J6").Value = 0.0209748055045382
$ws.Range("M6").Value = 44.56267066666667
$ws.Range("N6").Value = 133.688012
$ws.Range("O6").Value = 0.186578549555334
$ws.Range("P6").Value = 0.1926453158454455
$ws.Range("Q6").Value = 0.8290290708591112
$ws.Range("R6").Value = 7.461261637732001
$ws.Range("S6").Value = 0.002661716520073356
$ws.Range("T6").Value = 0.00404069803121855

$ws.Range("I7").Value = 0.02617205512618514
$ws.Range("J7").Value = 0.03848005474923699
$ws.Range("M7").Value = 98.170451
$ws.Range("N7").Value = 294.511353
$ws.Range("O7").Value = 0.4110278868558457
$ws.Range("P7").Value = 0.4243928215400083
$ws.Range("Q7").Value = 3.35055749263
$ws.Range("R7").Value = 30.15501743367
$ws.Range("S7").Value = 0.01075744451319058
$ws.Range("T7").Value = 0.01633065900804269

$ws.Range("I8").Value = 0.02617205512618514
$ws.Range("J8").Value = 0.03848005474923699
$ws.Range("O8").Value = 0.1453748421312515
$ws.Range("P8").Value = 0.150101833491052
$ws.Range("S8").Value = 0.003804758382219575
$ws.Range("T8").Value = 0.005775926770696534

$ws.Range("I9").Value = 0.02617205512618514
$ws.Range("J9").Value = 0.03848005474923699
$ws.Range("M9").Value = 38.82199566666667
$ws.Range("N9").Value = 116.465987
$ws.Range("O9").Value = 0.1625430328561575
$ws.Range("P9").Value = 0.1678282630971171
$ws.Range("Q9").Value = 1.324994712103333
$ws.Range("R9").Value = 11.92495240893
$ws.Range("S9").Value = 0.004254085216288676
$ws.Range("T9").Value = 0.006458040752446417

$ws.Range("I10").Value = 0.02617205512618514
$ws.Range("J10").Value = 0.03848005474923699
$ws.Range("M10").Value = 22.5647
$ws.Range("N10").Value = 45.1294
$ws.Range("O10").Value = 0.09447568860141126
$ws.Range("P10").Value = 0.06503176602637677
$ws.Range("Q10").Value = 0.7701332110000001
$ws.Range("R10").Value = 4.620799266000001
$ws.Range("S10").Value = 0.002472622930160436
$ws.Range("T10").Value = 0.002502425917134548

$ws.Range("I11").Value = 0.02617205512618514
$ws.Range("J11").Value = 0.03848005474923699
$ws.Range("M11").Value = 44.56267066666667
$ws.Range("N11").Value = 133.688012
$ws.Range("O11").Value = 0.186578549555334
$ws.Range("P11").Value = 0.1926453158454455
$ws.Range("Q11").Value = 1.520923949853333
$ws.Range("R11").Value = 13.68831554868
$ws.Range("S11").Value = 0.004883144084325866
$ws.Range("T11").Value = 0.007413002300916797

$ws.Range("G12").Value = 1.251329
$ws.Range("H12").Value = 2.502658
$ws.Range("I12").Value = 0.9595620149133938
$ws.Range("J12").Value = 0.9405451397462248
$ws.Range("M12").Value = 98.170451
$ws.Range("N12").Value = 294.511353
$ws.Range("O12").Value = 0.4110278868558457
$ws.Range("P12").Value = 0.4243928215400083
$ws.Range("Q12").Value = 122.843532279379
$ws.Range("R12").Value = 737.061193676274
$ws.Range("S12").Value = 0.3944067472969898
$ws.Range("T12").Value = 0.3991606056426418

$ws.Range("G13").Value = 1.251329
$ws.Range("H13").Value = 2.502658
$ws.Range("I13").Value = 0.9595620149133938
$ws.Range("J13").Value = 0.9405451397462248
$ws.Range("O13").Value = 0.1453748421312515
$ws.Range("P13").Value = 0.150101833491052
$ws.Range("Q13").Value = 43.44804740273801
$ws.Range("R13").Value = 260.6882844164281
$ws.Range("S13").Value = 0.1394961764331802
$ws.Range("T13").Value = 0.141177549957006

$ws.Range("G14").Value = 1.251329
$ws.Range("H14").Value = 2.502658
$ws.Range("I14").Value = 0.9595620149133938
$ws.Range("J14").Value = 0.9405451397462248
$ws.Range("M14").Value = 38.82199566666667
$ws.Range("N14").Value = 116.465987
$ws.Range("O14").Value = 0.1625430328561575
$ws.Range("P14").Value = 0.1678282630971171
$ws.Range("Q14").Value = 48.57908901557434
$ws.Range("R14").Value = 291.474534093446
$ws.Range("S14").Value = 0.1559701201175884
$ws.Range("T14").Value = 0.1578500571680442

$ws.Range("G15").Value = 1.251329
$ws.Range("H15").Value = 2.502658
$ws.Range("I15").Value = 0.9595620149133938
$ws.Range("J15").Value = 0.9405451397462248
$ws.Range("M15").Value = 22.5647
$ws.Range("N15").Value = 45.1294
$ws.Range("O15").Value = 0.09447568860141126
$ws.Range("P15").Value = 0.06503176602637677
$ws.Range("Q15").Value = 28.2358634863
$ws.Range("R15").Value = 112.9434539452
$ws.Range("S15").Value = 0.09065528211470053
$ws.Range("T15").Value = 0.06116531146522233

$ws.Range("G16").Value = 1.251329
$ws.Range("H16").Value = 2.502658
$ws.Range("I16").Value = 0.9595620149133938
$ws.Range("J16").Value = 0.9405451397462248
$ws.Range("M16").Value = 44.56267066666667
$ws.Range("N16").Value = 133.688012
$ws.Range("O16").Value = 0.186578549555334
$ws.Range("P16").Value = 0.1926453158454455
$ws.Range("Q16").Value = 55.76256212264934
$ws.Range("R16").Value = 334.5753727358961
$ws.Range("S16").Value = 0.1790336889509347
$ws.Range("T16").Value = 0.1811916155133102
